# Add names to User Cards
$wb = $excel.ActiveWorkbook

$wsUser = $wb.Worksheets.Item("User")
$wsApp  = $wb.Worksheets.Item("Application")

# Row -> Name mapping (column A) for the "User" sheet, in the same order the
# values were originally entered (this controls shared-string allocation
# order). Row 15 and 26/27 are intentionally left blank, matching the
# original layout.
$entries = @(
    @{ Row = 12; Name = "Bankers" }
    @{ Row = 16; Name = "Office Workers" }
    @{ Row = 25; Name = "Home Cinema Owners" }
    @{ Row = 17; Name = "MOOC Students" }
    @{ Row = 14; Name = "University Lecturers" }
    @{ Row = 7;  Name = "Tele Workers" }
    @{ Row = 20; Name = "Market Researchers" }
    @{ Row = 6;  Name = "Online shoppers" }
    @{ Row = 9;  Name = "Free to Play Gamers" }
    @{ Row = 10; Name = "Workplace Gamers" }
    @{ Row = 13; Name = "Premium News Readers" }
    @{ Row = 21; Name = "Social Media Shills" }
    @{ Row = 22; Name = "Tele Working Consultants" }
    @{ Row = 23; Name = "Video Game Bloggers" }
    @{ Row = 24; Name = "Mash up Artists" }
    @{ Row = 18; Name = "Game Review Excessive" }
    @{ Row = 19; Name = "Creative Professionals" }
    @{ Row = 5;  Name = "Fun loving socialites" }
    @{ Row = 8;  Name = "University Researchers" }
    @{ Row = 11; Name = "Documentary Watchers" }
)

foreach ($entry in $entries) {
    $wsUser.Cells.Item($entry.Row, 1).Value = $entry.Name
}

# Update the remembered selections to match the post-edit state.
$wsApp.Range("C31").Select()
$wsUser.Range("A11").Select()
$wsUser.Activate()
